$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (rows 2-15) from 10129 to 10195
$ws.Range("D2:D15").Value = 10195

# Update the active cell selection to G19
$ws.Range("G19").Select()
